# The two attendance records that occupied rows 15 and 16 were reordered
# (row 15 now holds the earlier-dated match, row 16 the later one). Sort
# the A15:Y16 block by the "date" column (E) ascending so the two rows
# swap places in full.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A15:Y16")
$key = $ws.Range("E15:E16")

$rng.Sort($key, 1)
